$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I2").Value = "b"
$ws.Range("J2").Value = "Acknowledge (Backchannel)"
$ws.Range("I10").Value = "b"
$ws.Range("J10").Value = "Acknowledge (Backchannel)"
$ws.Range("I18").Value = "ba"
$ws.Range("J18").Value = "Appreciation"
$ws.Range("I19").Value = "b"
$ws.Range("J19").Value = "Acknowledge (Backchannel)"
$ws.Range("I20").Value = "ba"
$ws.Range("J20").Value = "Appreciation"
$ws.Range("I30").Value = "%"
$ws.Range("J30").Value = "Uninterpretable"
$ws.Range("I49").Value = "sv"
$ws.Range("J49").Value = "Statement-opinion"
$ws.Range("I51").Value = "b"
$ws.Range("J51").Value = "Acknowledge (Backchannel)"
$ws.Range("I53").Value = "aa"
$ws.Range("J53").Value = "Agree/Accept"
$ws.Range("I55").Value = "b"
$ws.Range("J55").Value = "Acknowledge (Backchannel)"
$ws.Range("I57").Value = "sv"
$ws.Range("J57").Value = "Statement-opinion"
$ws.Range("I59").Value = "%"
$ws.Range("J59").Value = "Uninterpretable"
$ws.Range("I60").Value = "b"
$ws.Range("J60").Value = "Acknowledge (Backchannel)"
$ws.Range("I65").Value = "aa"
$ws.Range("J65").Value = "Agree/Accept"
$ws.Range("I72").Value = "sv"
$ws.Range("J72").Value = "Statement-opinion"
$ws.Range("I102").Value = "sd"
$ws.Range("J102").Value = "Statement-non-opinion"
$ws.Range("I114").Value = "b"
$ws.Range("J114").Value = "Acknowledge (Backchannel)"
$ws.Range("I118").Value = "b"
$ws.Range("J118").Value = "Acknowledge (Backchannel)"
$ws.Range("I132").Value = "ba"
$ws.Range("J132").Value = "Appreciation"
$ws.Range("I133").Value = "sd"
$ws.Range("J133").Value = "Statement-non-opinion"
$ws.Range("I140").Value = "ba"
$ws.Range("J140").Value = "Appreciation"
$ws.Range("I142").Value = "aa"
$ws.Range("J142").Value = "Agree/Accept"
$ws.Range("I148").Value = "b"
$ws.Range("J148").Value = "Acknowledge (Backchannel)"
$ws.Range("I169").Value = "aa"
$ws.Range("J169").Value = "Agree/Accept"
$ws.Range("I177").Value = "aa"
$ws.Range("J177").Value = "Agree/Accept"
$ws.Range("I184").Value = "sd"
$ws.Range("J184").Value = "Statement-non-opinion"
$ws.Range("I188").Value = "sd"
$ws.Range("J188").Value = "Statement-non-opinion"
$ws.Range("I197").Value = "%"
$ws.Range("J197").Value = "Uninterpretable"
$ws.Range("I200").Value = "b"
$ws.Range("J200").Value = "Acknowledge (Backchannel)"
$ws.Range("I216").Value = "sd"
$ws.Range("J216").Value = "Statement-non-opinion"
$ws.Range("I239").Value = "sd"
$ws.Range("J239").Value = "Statement-non-opinion"
$ws.Range("I246").Value = "aa"
$ws.Range("J246").Value = "Agree/Accept"
$ws.Range("I276").Value = "b"
$ws.Range("J276").Value = "Acknowledge (Backchannel)"
$ws.Range("I284").Value = "aa"
$ws.Range("J284").Value = "Agree/Accept"
$ws.Range("I285").Value = "aa"
$ws.Range("J285").Value = "Agree/Accept"
$ws.Range("I297").Value = "ba"
$ws.Range("J297").Value = "Appreciation"
$ws.Range("I298").Value = "b"
$ws.Range("J298").Value = "Acknowledge (Backchannel)"
$ws.Range("I304").Value = "ba"
$ws.Range("J304").Value = "Appreciation"
$ws.Range("I305").Value = "b"
$ws.Range("J305").Value = "Acknowledge (Backchannel)"
$ws.Range("I312").Value = "ba"
$ws.Range("J312").Value = "Appreciation"
